# Fruta / hortaliza, semanal
# Update the weekly price data rows (2-12) on the sheet to reflect the
# new period's values (dates, quality grade, volume, prices, origin region).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: D=Fecha, L=Calidad, M=Volumen, N=Precio minimo,
# O=Precio maximo, P=Precio promedio ponderado, R=Origen, S=Precio $/Kg
$rows = @(
    @{ Row = 2;  D = 44295; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región Metropolitana";   S = 506 },
    @{ Row = 3;  D = 44295; L = "Segunda";  M = 16; N = 195000; O = 200000; P = 197500; R = "Región Metropolitana";   S = 439 },
    @{ Row = 4;  D = 44309; L = "Especial"; M = 20; N = 305000; O = 310000; P = 307500; R = "Provincia de Cachapoal"; S = 683 },
    @{ Row = 5;  D = 44309; L = "Primera";  M = 20; N = 285000; O = 290000; P = 287500; R = "Provincia de Cachapoal"; S = 639 },
    @{ Row = 6;  D = 44309; L = "Segunda";  M = 20; N = 255000; O = 260000; P = 257500; R = "Provincia de Cachapoal"; S = 572 },
    @{ Row = 7;  D = 44316; L = "Especial"; M = 20; N = 255000; O = 260000; P = 257500; R = "Región de O'Higgins";    S = 572 },
    @{ Row = 8;  D = 44316; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región de O'Higgins";    S = 506 },
    @{ Row = 9;  D = 44273; L = "Especial"; M = 10; N = 255000; O = 260000; P = 257500; R = "Región de O'Higgins";    S = 572 },
    @{ Row = 10; D = 44273; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región de O'Higgins";    S = 506 },
    @{ Row = 11; D = 44294; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región Metropolitana";   S = 506 },
    @{ Row = 12; D = 44294; L = "Segunda";  M = 16; N = 195000; O = 200000; P = 197500; R = "Región Metropolitana";   S = 439 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $item.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $item.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $item.R   # R - Origen
    $ws.Cells.Item($r, 19).Value = $item.S   # S - Precio $/Kg
}
